# Running all the suites
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Update the test result for row 12 (ProfileFollowerTest):
# Runmode column (C) -> "Y" (was "N")
# Result column (D) -> "SKIP" (was "PASS")
$ws.Range("C12").Value = "Y"
$ws.Range("D12").Value = "SKIP"

# Move the active selection to C12
$ws.Activate()
$ws.Range("C12").Select()
